# Generate Report for Handoff
# Rebuild Overview / zh-cn / de-de sheets with the new handoff report data.

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276   # BGR int for RGB 6495ED (matches existing HyperLink font color)
$DATE_FMT = "yyyy-mm-dd HH:mm:ss"

function Set-HyperlinkCell($ws, $row, $col, $displayText, $url) {
    $cell = $ws.Cells.Item($row, $col)
    $ws.Hyperlinks.Add($cell, $url, "", "", $displayText)
    $cell.Font.Underline = $true
    $cell.Font.Color = $HYPERLINK_COLOR
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Clear()
$ws1.Hyperlinks.Delete()

$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "zh-cn"
$ws1.Range("C1").Value = "de-de"
$ws1.Range("D1").Value = "Latest Handoff Date"

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-24 03:06:20"

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-24 03:06:20"

$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-24 03:06:20"

Set-HyperlinkCell $ws1 2 1 "811e9129-2027-4d9c-8b7d-594bbd018693.md" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/811e9129-2027-4d9c-8b7d-594bbd018693.md"
Set-HyperlinkCell $ws1 3 1 "88f2dff0-dd20-4021-9189-d32f8bd30135.png" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/88f2dff0-dd20-4021-9189-d32f8bd30135.png"
Set-HyperlinkCell $ws1 4 1 "da01fbcb-273b-4d5c-b732-88ef238269a4.png" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/da01fbcb-273b-4d5c-b732-88ef238269a4.png"

$ws1.Range("D2:D4").NumberFormat = $DATE_FMT

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Clear()
$ws2.Hyperlinks.Delete()

$headers = @("Source File Name","File Extension","Status","Latest Handoff File","Latest Handoff Datetime","Latest Target File","Latest Handback File","Latest Handback DateTime","Reference Tokens","Handoff Reason","Dependency From","Error Detail")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - md file
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-24 03:06:12"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("J2").Value = "Include"

# Row 3 - first png (dependency)
$ws2.Range("B3").Value = ".png"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "2016-03-24 03:06:12"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "IsDependency"
$ws2.Range("K3").Value = "e2e\811e9129-2027-4d9c-8b7d-594bbd018693.md"

# Row 4 - second png (dependency)
$ws2.Range("B4").Value = ".png"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("E4").Value = "2016-03-24 03:06:12"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "IsDependency"
$ws2.Range("K4").Value = "e2e\811e9129-2027-4d9c-8b7d-594bbd018693.md"

Set-HyperlinkCell $ws2 2 1 "811e9129-2027-4d9c-8b7d-594bbd018693.md" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/811e9129-2027-4d9c-8b7d-594bbd018693.md"
Set-HyperlinkCell $ws2 2 4 "811e9129-2027-4d9c-8b7d-594bbd018693.ccfdc9013870b5f0f0cc8648963c2f5737667d9e.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/811e9129-2027-4d9c-8b7d-594bbd018693.ccfdc9013870b5f0f0cc8648963c2f5737667d9e.zh-cn.xlf"
Set-HyperlinkCell $ws2 3 1 "88f2dff0-dd20-4021-9189-d32f8bd30135.png" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/88f2dff0-dd20-4021-9189-d32f8bd30135.png"
Set-HyperlinkCell $ws2 3 4 "e9f82d5a7f6850c56d286ff8f502bd51fabf3f7f.png" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/e9f82d5a7f6850c56d286ff8f502bd51fabf3f7f.png"
Set-HyperlinkCell $ws2 4 1 "da01fbcb-273b-4d5c-b732-88ef238269a4.png" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/da01fbcb-273b-4d5c-b732-88ef238269a4.png"
Set-HyperlinkCell $ws2 4 4 "4a4ac5b595002c4e4d6455813836a50cc661faeb.png" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/4a4ac5b595002c4e4d6455813836a50cc661faeb.png"

$ws2.Range("E2:E4").NumberFormat = $DATE_FMT
$ws2.Range("H2:H4").NumberFormat = $DATE_FMT

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Clear()
$ws3.Hyperlinks.Delete()

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - md file
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-24 03:06:20"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("J2").Value = "Include"

# Row 3 - first png (dependency)
$ws3.Range("B3").Value = ".png"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "2016-03-24 03:06:20"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "IsDependency"
$ws3.Range("K3").Value = "e2e\811e9129-2027-4d9c-8b7d-594bbd018693.md"

# Row 4 - second png (dependency)
$ws3.Range("B4").Value = ".png"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("E4").Value = "2016-03-24 03:06:20"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "IsDependency"
$ws3.Range("K4").Value = "e2e\811e9129-2027-4d9c-8b7d-594bbd018693.md"

Set-HyperlinkCell $ws3 2 1 "811e9129-2027-4d9c-8b7d-594bbd018693.md" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/811e9129-2027-4d9c-8b7d-594bbd018693.md"
Set-HyperlinkCell $ws3 2 4 "811e9129-2027-4d9c-8b7d-594bbd018693.ccfdc9013870b5f0f0cc8648963c2f5737667d9e.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/811e9129-2027-4d9c-8b7d-594bbd018693.ccfdc9013870b5f0f0cc8648963c2f5737667d9e.de-de.xlf"
Set-HyperlinkCell $ws3 3 1 "88f2dff0-dd20-4021-9189-d32f8bd30135.png" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/88f2dff0-dd20-4021-9189-d32f8bd30135.png"
Set-HyperlinkCell $ws3 3 4 "e9f82d5a7f6850c56d286ff8f502bd51fabf3f7f.png" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/e9f82d5a7f6850c56d286ff8f502bd51fabf3f7f.png"
Set-HyperlinkCell $ws3 4 1 "da01fbcb-273b-4d5c-b732-88ef238269a4.png" "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/da01fbcb-273b-4d5c-b732-88ef238269a4.png"
Set-HyperlinkCell $ws3 4 4 "4a4ac5b595002c4e4d6455813836a50cc661faeb.png" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/4a4ac5b595002c4e4d6455813836a50cc661faeb.png"

$ws3.Range("E2:E4").NumberFormat = $DATE_FMT
$ws3.Range("H2:H4").NumberFormat = $DATE_FMT
